$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 gets translated (localized) to a new album entry.
$ws.Range("A7").Value = "Kokoroko"
$ws.Range("B7").Value = "Could We Be More"
$ws.Range("C7").Value = 44774
$ws.Range("D7").Value = 4.7
$ws.Range("E7").Value = 30

# Move the active selection to C8, matching the author's last click.
$ws.Range("C8").Select() | Out-Null
